$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.26"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.68"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.478"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05699"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.370"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8018"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.039"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1520"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07391"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03154"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03001"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09300"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001644"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.421"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04696"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005868"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006352"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005054"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001043"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0003134"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.776"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.433"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.124"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3286"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006970"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1047"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002970"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008573"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005825"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005498"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6823"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.009310"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
